$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Créditos-aula: 4 -> 2  (row 5). Force text so it stays a literal string
# like the original ("4") instead of being auto-coerced to a number.
$ws.Range("B5:C5").NumberFormat = "@"
$ws.Range("B5").Value = "2"
$ws.Range("C5").Value = "2"

# Carga horária: 60 h -> 30 h (row 7) - already non-numeric text, no coercion risk.
$ws.Range("B7").Value = "30 h"
$ws.Range("C7").Value = "30 h"

# Ativação: 01/01/2020 -> 01/01/2022 (row 8). Force text so the date string
# stays literal text instead of being parsed into a date serial value.
$ws.Range("B8:C8").NumberFormat = "@"
$ws.Range("B8").Value = "01/01/2022"
$ws.Range("C8").Value = "01/01/2022"

# Programa resumido: (row 15)
$programaResumido = "Escoamento permanente uniforme em condutos forçados e perdas de carga nas tubulações. Instalações de bombeamento e bombas hidráulicas. Escoamento permanente uniforme em condutos livres, resistência ao escoamento e perdas de carga nos canais naturais e artificiais. Estudo da carga específica em canais."
$ws.Range("B15").Value = $programaResumido
$ws.Range("C15").Value = $programaResumido

# Short syllabus: (row 16)
$shortSyllabus = "Permanent flow in pressure conduits. Head losses in pipelines. Pumping stations and hydraulic pumps. Permanente flow in free surface conduits, head losses in natural and artificial channels. Specific Energy in channels."
$ws.Range("B16").Value = $shortSyllabus
$ws.Range("C16").Value = $shortSyllabus

# Programa: (row 17)
$programa = "- Hidrostática,- piezometria,- conservação da massa e quantidade de movimento,- Escoamentos Permanentes em Condutos Forçados,- Resistência ao Escoamento e Perdas de Carga,- Bombas e sistemas de recalque,- Escoamento Permanente Uniforme em Condutos Livres,- Resistência ao escoamento e Perdas de Carga,- Canais regulares e naturais,- Carga Específica,- Escoamento Permanente Gradualmente Variado,- Cálculo da linha d’água,- Ressalto Hidráulico."
$ws.Range("B17").Value = $programa
$ws.Range("C17").Value = $programa

# Syllabus: (row 18)
$syllabus = "- Hydrostatic,- piezometry,- conservation of mass and amount of movement,- Permanent flows in Pipes,- Resistance to Flow and Load Losses,- Pumps and discharge systems,- Uniform Permanent Flow in Free Flues,- Resistance to flow and head losses,- Regular and natural channels,- Specific Load,- Gradually Variable Permanent Flow,- Calculation of the water line,- Hydraulic boss."
$ws.Range("B18").Value = $syllabus
$ws.Range("C18").Value = $syllabus
